$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record needs to be inserted as the first (most recent) entry
# of the Ciboulette price table, which is sorted with the newest date at
# the bottom of the existing range (row 163) and dates shifting down as
# older records move further down. Insert a new row at row 163, pushing
# rows 163:176 down to 164:177, then fill the new row with the same
# constant attributes as the rest of the block plus the new date/volume.

$ws.Rows("163:163").Insert()

$newRow = 163

$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44461
$ws.Cells.Item($newRow, 4).Style = $ws.Cells.Item($newRow + 1, 4).Style
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = 100112039
$ws.Cells.Item($newRow, 7).Value = "Ciboulette"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 160
$ws.Cells.Item($newRow, 11).Value = 1500
$ws.Cells.Item($newRow, 12).Value = 1500
$ws.Cells.Item($newRow, 13).Value = 1500
$ws.Cells.Item($newRow, 14).Value = "$/docena de atados"
$ws.Cells.Item($newRow, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($newRow, 16).Value = 500
$ws.Cells.Item($newRow, 17).Value = 3
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
